# Updated drug catalog specification page
# - Refresh the static "last updated" date placeholder text (slide master,
#   notes master, and every slide layout) from 17/06/2021 to 13/08/2021.
# - Rename the "DrugIngredient profile of Ingredient" label to
#   "DrugSubstance profile of Substance" on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "17/06/2021"
$newDate = "13/08/2021"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -or $shp.Name -like "*eserv*date*") {
            if ($shp.HasTextFrame) {
                if ($shp.TextFrame.HasText) {
                    if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                        $shp.TextFrame.TextRange.Text = $newDate
                    }
                }
            }
        }
    }
}

# NotesMaster's placeholder shapes don't accept direct TextRange writes in
# this host, but its HeadersFooters.DateAndTime.Text setter does persist
# (its getter is a stub, so we just assign unconditionally - the target
# deck only ever carries the one fixed date in that slot).
function Set-HeaderFooterDate($container) {
    $dt = $container.HeadersFooters.DateAndTime
    $dt.Text = $newDate
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout attached to the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master
$notesMaster = $p.NotesMaster
Set-HeaderFooterDate $notesMaster

# Slide 1: rename the DrugIngredient label to DrugSubstance
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "DrugIngredient profile of Ingredient") {
                $shp.TextFrame.TextRange.Text = "DrugSubstance profile of Substance"
            }
        }
    }
}
